$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = [double]"337.2695627278077"
$ws.Range("F2").Value = [double]"175.0122316609742"
$ws.Range("G2").Value = [double]"-8.58635347429663e-11"
$ws.Range("K2").Value = [double]"-0.4525505391315189"
$ws.Range("C3").Value = [double]"625.8612525615872"
$ws.Range("F3").Value = [double]"397.1059129019827"
$ws.Range("G3").Value = [double]"245.6123408855054"
$ws.Range("K3").Value = [double]"-1.024511054103239"
$ws.Range("C4").Value = [double]"822.9134908163769"
$ws.Range("F4").Value = [double]"648.3402844044327"
$ws.Range("G4").Value = [double]"802.5533206762635"
$ws.Range("K4").Value = [double]"-1.673079553805425"
$ws.Range("C5").Value = [double]"897.9960599337786"
$ws.Range("F5").Value = [double]"897.763784803767"
$ws.Range("G5").Value = [double]"1714.343376201689"
$ws.Range("K5").Value = [double]"-2.313735200304391"
$ws.Range("C6").Value = [double]"926.9347789622033"
$ws.Range("F6").Value = [double]"1156.953952186595"
$ws.Range("G6").Value = [double]"2978.324550276093"
$ws.Range("K6").Value = [double]"-2.976464721003201"
$ws.Range("C7").Value = [double]"935.5258614050725"
$ws.Range("F7").Value = [double]"1425.831136761117"
$ws.Range("G7").Value = [double]"4607.225239504814"
$ws.Range("K7").Value = [double]"-3.662948383454637"
$ws.Range("C8").Value = [double]"930.9463650054442"
$ws.Range("F8").Value = [double]"1701.377148470918"
$ws.Range("G8").Value = [double]"6614.683984933854"
$ws.Range("K8").Value = [double]"-4.365648273418756"
$ws.Range("C9").Value = [double]"872.8310726118056"
$ws.Range("F9").Value = [double]"1987.147194107323"
$ws.Range("G9").Value = [double]"9007.839350411174"
$ws.Range("K9").Value = [double]"-4.876124632317584"
$ws.Range("C10").Value = [double]"903.6782435345085"
$ws.Range("F10").Value = [double]"2297.737444150939"
$ws.Range("G10").Value = [double]"11649.68579295754"
$ws.Range("K10").Value = [double]"-6.249526756895477"
$ws.Range("C11").Value = [double]"853.1687190745885"
$ws.Range("F11").Value = [double]"2614.389296533142"
$ws.Range("G11").Value = [double]"14983.8166790652"
$ws.Range("K11").Value = [double]"-6.64920211999027"
$ws.Range("C12").Value = [double]"853.4171916819225"
$ws.Range("F12").Value = [double]"-2614.439928077682"
$ws.Range("G12").Value = [double]"18646.69199200964"
$ws.Range("K12").Value = [double]"6.651181569588934"
$ws.Range("C13").Value = [double]"904.2670374912242"
$ws.Range("G13").Value = [double]"15077.63455762254"
$ws.Range("K13").Value = [double]"6.25367608101158"
$ws.Range("C14").Value = [double]"873.5564153245809"
$ws.Range("F14").Value = [double]"-1987.510904179602"
$ws.Range("G14").Value = [double]"11678.79904100106"
$ws.Range("K14").Value = [double]"4.88023692952197"
$ws.Range("C15").Value = [double]"931.6159064490297"
$ws.Range("F15").Value = [double]"-1701.676930506009"
$ws.Range("G15").Value = [double]"9012.017504123072"
$ws.Range("K15").Value = [double]"4.368844526648956"
$ws.Range("C16").Value = [double]"936.2359146907713"
$ws.Range("F16").Value = [double]"-1426.138684620535"
$ws.Range("G16").Value = [double]"6616.189777221574"
$ws.Range("K16").Value = [double]"3.665775890114592"
$ws.Range("C17").Value = [double]"927.6396865495901"
$ws.Range("F17").Value = [double]"-1157.204843849326"
$ws.Range("G17").Value = [double]"4608.298301604259"
$ws.Range("K17").Value = [double]"2.978766731585818"
$ws.Range("C18").Value = [double]"898.6924593665523"
$ws.Range("F18").Value = [double]"-897.9716041544795"
$ws.Range("G18").Value = [double]"2979.044593351941"
$ws.Range("K18").Value = [double]"2.31555943445677"
$ws.Range("C19").Value = [double]"823.5477617979968"
$ws.Range("F19").Value = [double]"-648.4794598771598"
$ws.Range("G19").Value = [double]"1715.027634685314"
$ws.Range("K19").Value = [double]"1.674390738614533"
$ws.Range("C20").Value = [double]"626.3769706185988"
$ws.Range("G20").Value = [double]"803.1070853603685"
$ws.Range("K20").Value = [double]"1.025368515122934"
$ws.Range("C21").Value = [double]"337.5478961334583"
$ws.Range("F21").Value = [double]"-175.0694776766458"
$ws.Range("G21").Value = [double]"245.6809126370629"
$ws.Range("K21").Value = [double]"0.4529288330461948"
